# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (the stock Office palette),
#     used only by the notes master.
#   ppt/theme/theme2.xml -> "Integral" (a custom green/teal palette),
#     used by the slide master / the presentation as a whole.
#
# The authored edit swaps the two themes' contents: the deck-wide theme
# (theme2.xml) becomes the plain "Office Theme" colour palette (and,
# conversely, the notes-only theme becomes "Integral" — but the notes
# master's theme part is not independently addressable through the
# PowerPoint object model, only the presentation's single active theme
# colour scheme is). We reproduce the reachable half of that swap: the
# presentation's live ThemeColorScheme (backed by ppt/theme/theme2.xml)
# is repointed from the "Integral" palette to the stock Office palette,
# colour-slot by colour-slot, via ThemeColorScheme.Colors(i).RGB — the
# documented, file-system-free way to edit a theme's colours in this
# host.

function ComRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme, in ThemeColorSchemeIndex order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $tcs.Colors($i).RGB = ComRGB $officeThemeColors[$i - 1]
}
